$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay plain text (these look numeric but are
# formatted strings like "591.94" or "1.00" with significant trailing
# zeros / thousands-dot notation that must survive verbatim).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.320.01"
$ws.Range("E2").Value = "  -3.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.773.61"
$ws.Range("E3").Value = "  +1.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "591.94"
$ws.Range("E5").Value = "  -3.55%  "

# Row 6 - Solana
$ws.Range("D6").Value = "171.74"
$ws.Range("E6").Value = "  -4.00%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.779.48"
$ws.Range("E7").Value = "  +1.31%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -1.83%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -4.35%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  -4.71%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -4.17%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "37.72"
$ws.Range("E13").Value = "  -5.07%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -3.93%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.397.73"
$ws.Range("E15").Value = "  +1.09%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.765.18"
$ws.Range("E16").Value = "  +0.92%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.392.77"
$ws.Range("E17").Value = "  -3.34%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -4.72%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "7.10"
$ws.Range("E19").Value = "  -5.16%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "16.12"
$ws.Range("E20").Value = "  -1.20%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "486.72"
$ws.Range("E21").Value = "  -2.93%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "9.14"
$ws.Range("E22").Value = "  -0.04%  "

# Row 23 - Polygon
$ws.Range("D23").Value = "0.721"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "84.02"
$ws.Range("E24").Value = "  -2.45%  "

# Row 25 - Fetch.AI
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  -9.43%  "

# Row 26 - PEPE
$ws.Range("D26").Value = "0.0000138"
$ws.Range("E26").Value = "  +1.60%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "12.21"
$ws.Range("E27").Value = "  -5.36%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  -10.26%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.03%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "2.91"
$ws.Range("E30").Value = "  -0.30%  "

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -2.40%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "32.30"
$ws.Range("E32").Value = "  +6.55%  "

# Row 33 - NEARProtocol
$ws.Range("D33").Value = "7.70"
$ws.Range("E33").Value = "  -4.32%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -4.77%  "

# Row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.25%  "

# Row 36 - Mantle
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -3.93%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  -1.57%  "

# Row 38 - Filecoin
$ws.Range("D38").Value = "5.72"
$ws.Range("E38").Value = "  -6.25%  "

# Row 39 - Bittensor
$ws.Range("D39").Value = "451.73"
$ws.Range("E39").Value = "  +3.33%  "

# Row 40 - TheGraph
$ws.Range("D40").Value = "0.323"
$ws.Range("E40").Value = "  -7.45%  "

# Row 41 - OKB
$ws.Range("D41").Value = "48.79"
$ws.Range("E41").Value = "  -1.74%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -3.92%  "

# Row 43 - dogwifhat
$ws.Range("D43").Value = "2.85"
$ws.Range("E43").Value = "  -6.18%  "

# Row 44 - Cosmos
$ws.Range("D44").Value = "8.24"
$ws.Range("E44").Value = "  -3.78%  "

# Row 45 - Arweave
$ws.Range("D45").Value = "41.16"
$ws.Range("E45").Value = "  -9.97%  "

# Row 46 - Maker
$ws.Range("D46").Value = "2.808.89"
$ws.Range("E46").Value = "  -4.87%  "

# Row 47 - was USDe, now Monero
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "139.30"
$ws.Range("E47").Value = "  +0.18%  "

# Row 48 - was Monero, now USDe
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "1.00"

# Row 49 - VeChain
$ws.Range("D49").Value = "0.0347"
$ws.Range("E49").Value = "  -3.73%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "25.83"
$ws.Range("E50").Value = "  -4.84%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "23.15"
$ws.Range("E51").Value = "  +8.67%  "
